$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Skeleton knight hp/ht value changes from 0/10 to 0/12
$ws.Range("C8").Value = "0/12"

# New loot entry for skeleton knight: humanity
$ws.Range("J8").Value = "humanity"

# Move active selection to J14 (new target cell after edit)
$ws.Range("J14").Select()
